$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.58250000000003
$ws.Range("E4").Value = 14.21209999999999

$ws.Range("E5").Value = 13.3471

$ws.Range("A7").Value = -21.55880000000001

$ws.Range("E8").Value = 14.38649999999999

$ws.Range("A16").Value = -20.06369999999999
$ws.Range("E16").Value = 13.09750000000001
